$d = $word.ActiveDocument

# --- Programa (PT) paragraph: split "... estatística. " / "2. Introdução ..." with a line break ---
$d.Content.Find.Execute(
    "estatística. 2. Introdução",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "estatística. ^l2. Introdução",
    2
)

# --- Programa (EN, italic) paragraph: split "... analysis." / "2. Introduction ..." with a line break ---
$d.Content.Find.Execute(
    "statistical analysis.2. Introduction",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "statistical analysis.^l2. Introduction",
    2
)

# --- Avaliação / Critério run: break up "MF = ... Onde: ... P1 ... P2 ... MF ..." into separate lines ---
$d.Content.Find.Execute(
    "(P1 + P2)/2 Onde: P1",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "(P1 + P2)/2 ^lOnde: P1",
    2
)

$d.Content.Find.Execute(
    "Onde: P1 é a nota obtida pela avaliação da resolução de problemas referentes ao item 1 do Programa do curso;P2",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Onde: ^lP1 é a nota obtida pela avaliação da resolução de problemas referentes ao item 1 do Programa do curso;^lP2",
    2
)

$d.Content.Find.Execute(
    "item 2 do Programa do curso;MF é a média final",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "item 2 do Programa do curso;^lMF é a média final",
    2
)

# --- Bibliografia run: break each numbered reference onto its own line ---
$d.Content.Find.Execute(
    "1968.2) FINLAYSON",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "1968.^l2) FINLAYSON",
    2
)

$d.Content.Find.Execute(
    "McGraw-Hill, 19803) CONSTANTINIDES",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "McGraw-Hill, 1980^l3) CONSTANTINIDES",
    2
)

$d.Content.Find.Execute(
    "1999.4) CUTLIP",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "1999.^l4) CUTLIP",
    2
)

$d.Content.Find.Execute(
    "2008.5) DYMENT",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "2008.^l5) DYMENT",
    2
)

$d.Content.Find.Execute(
    "2015.6) FARES",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "2015.^l6) FARES",
    2
)

Write-Host "Edits applied"
